$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("HG").Insert()

$ws.Range("HG2").Value = "vctrs"
$ws.Range("HG3").Value = "list_drop_empty"

for ($r = 4; $r -le 60; $r++) {
    if ($r -eq 22) {
        $ws.Range("HG$r").Value = 1
    } else {
        $ws.Range("HG$r").Value = 0
    }
}

Write-Host "Done"
